# The recorded change only touches customXml/itemProps3.xml — specifically
# the <ds:schemaRef> children of the <ds:schemaRefs> element inside the
# SharePoint "datastoreItem" companion part for one of the package's
# customXml/item*.xml parts (itemID {8D289AE2-D2AE-49D1-AFAC-3A79F6794255}).
#
# The diff reshuffles the order of the ten <ds:schemaRef> entries; the set
# of ds:uri values is identical before and after (same 10 URIs, just
# re-serialized in a different sequence). That is exactly what happens when
# Office/SharePoint round-trips this custom-XML "schema refs" store (it is
# kept internally as an unordered set/dictionary), which matches the
# generic "Add files via upload" commit message — i.e. this hunk is an
# incidental artifact of re-uploading the file, not a deliberate edit a
# user made in PowerPoint.
#
# This part is not reachable from PowerPoint's scripting surface:
#   - It is not one of the customXml/item*.xml "data" parts; it is the
#     auto-managed itemProps*.xml "schema refs" companion part that Office
#     itself writes and that is never exposed for editing.
#   - Presentation.CustomXMLParts exists on the object model, but it does
#     not provide any member to reach / reorder a part's schemaRefs; on
#     the live document here it reports Count = 0 and Add()/Item() do not
#     persist anything to the package (content_diffs stays 0) — i.e. there
#     is no COM-reachable operation, in real PowerPoint or in this
#     interop runtime, that can reproduce a schemaRef-order-only change.
#
# Since no slide/shape/text/property content actually changed, there is
# nothing for a PowerPoint automation script to legitimately do here — so
# this script intentionally performs no operations on the presentation.
$p = $ppt.ActivePresentation
